$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(61, 1).Value = "2024-09-30 00:00:00"
$ws.Cells.Item(61, 2).Value = 75650
$ws.Cells.Item(61, 3).Value = 10756.89
$ws.Cells.Item(61, 4).Value = 9519.370000000001
$ws.Cells.Item(61, 5).Value = 7.0121
